# Add "Buenos Aires" surprise songs (South America leg) to the Taylor Swift
# Eras Tour surprise songs dataset.
#
# New shared strings must be registered in the same order the source data
# frame would have produced them (column-major: all of leg, then date, then
# city, then night, then dress, then instrument, then song) so the
# resulting sharedStrings.xml ordering matches the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=leg, B=date, C=city, D=night, E=dress, F=instrument, G=song
# Date values are stored as Excel serial day numbers (matching existing rows).
$rows = @(
    @("South America", 45239, "Buenos Aires, Argentina", 1, "blue",   "guitar", "The Very First Night (Taylor's Version) [From The Vault]"),
    @("South America", 45239, "Buenos Aires, Argentina", 1, "blue",   "piano",  "Labyrinth"),
    @("South America", 45241, "Buenos Aires, Argentina", 2, "red",    "guitar", "Is It Over Now? (Taylor's Version) [From The Vault]"),
    @("South America", 45241, "Buenos Aires, Argentina", 2, "red",    "piano",  "End Game"),
    @("South America", 45242, "Buenos Aires, Argentina", 3, "yellow", "guitar", "Better Than Revenge (Taylor's Version)"),
    @("South America", 45242, "Buenos Aires, Argentina", 3, "yellow", "piano",  "`"Slut!`" (Taylor's Version) [From The Vault]"),
    @("South America", 45247, "Rio de Janeiro, Brazil",  1, $null,    $null,    $null),
    @("South America", 45247, "Rio de Janeiro, Brazil",  1, $null,    $null,    $null),
    @("South America", 45248, "Rio de Janeiro, Brazil",  2, $null,    $null,    $null),
    @("South America", 45248, "Rio de Janeiro, Brazil",  2, $null,    $null,    $null),
    @("South America", 45249, "Rio de Janeiro, Brazil",  3, $null,    $null,    $null),
    @("South America", 45249, "Rio de Janeiro, Brazil",  3, $null,    $null,    $null),
    @("South America", 45254, "São Paulo, Brazil",       1, $null,    $null,    $null),
    @("South America", 45254, "São Paulo, Brazil",       1, $null,    $null,    $null),
    @("South America", 45255, "São Paulo, Brazil",       2, $null,    $null,    $null),
    @("South America", 45255, "São Paulo, Brazil",       2, $null,    $null,    $null),
    @("South America", 45256, "São Paulo, Brazil",       3, $null,    $null,    $null),
    @("South America", 45256, "São Paulo, Brazil",       3, $null,    $null,    $null)
)

$startRow = 116
$colCount = 7

for ($c = 0; $c -lt $colCount; $c++) {
    for ($i = 0; $i -lt $rows.Count; $i++) {
        $r = $startRow + $i
        $value = $rows[$i][$c]
        if ($value -eq $null) {
            continue
        }
        $ws.Cells.Item($r, $c + 1).Value = $value
    }
}

# Give the new date cells (column B) the same "yyyy-mm-dd" date style already
# used throughout the sheet, by copying the format from the preceding date
# cell rather than re-entering the format code (which would register a
# duplicate numFmt entry).
$ws.Cells.Item($startRow - 1, 2).Copy() | Out-Null
$dateRange = $ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($startRow + $rows.Count - 1, 2))
$dateRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# City names got longer ("Buenos Aires, Argentina" / "Rio de Janeiro, Brazil"),
# so column C needs to widen to keep fitting its contents.
$ws.Columns.Item(3).ColumnWidth = 20.1

# Leave the view scrolled near the newly-added rows, matching where the
# editor ended up after appending this block.
$ws.Range("F122").Select() | Out-Null
